$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write in column-major order (B2,B3,B4,C2,C3,C4,D2,D3,D4) so that the
# newly created shared-string entries are appended in the same order as
# the original strings they replace.
$ws.Range("B2").Value = "'0.17"
$ws.Range("B3").Value = "'-0.01"
$ws.Range("B4").Value = "'-0.09"

$ws.Range("C2").Value = "44.29***"
$ws.Range("C3").Value = "2.21***"
$ws.Range("C4").Value = "'0.98"

$ws.Range("D2").Value = "'-0.89"
$ws.Range("D3").Value = "0.46***"
$ws.Range("D4").Value = "0.82*"
